# Daily Report update: append 2026-01-06 (serial 46028) block
# mirrors the prior day's (2026-01-05 / serial 46027) depository rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

$newDate = 46028
$startRow = 46

# Each entry: Region_Type label, PREV_TOTAL, RECEIVED, WITHDRAWN, NET_CHANGE, ADJUSTMENT, TOTAL_TODAY
$rows = @(
    @('ASAHI DEPOSITORY LLC Registered', 0, 0, 0, 0, 0, 0),
    @('ASAHI DEPOSITORY LLC Eligible', 0, 0, 0, 0, 0, 0),
    @('BRINK''S, INC. Registered', 90027.72500000001, 0, 0, 0, 0, 90027.72500000001),
    @('BRINK''S, INC. Eligible', 5744.711, 0, 0, 0, 0, 5744.711),
    @('CNT DEPOSITORY, INC. Registered', 1246.06, 0, 0, 0, 0, 1246.06),
    @('CNT DEPOSITORY, INC. Eligible', 0, 0, 0, 0, 0, 0),
    @('DELAWARE DEPOSITORY Registered', 1633.941, 0, 0, 0, 0, 1633.941),
    @('DELAWARE DEPOSITORY Eligible', 18509.729, 0, 0, 0, 0, 18509.729),
    @('HSBC BANK, USA Registered', 1295.223, 0, 0, 0, 0, 1295.223),
    @('HSBC BANK, USA Eligible', 9281.978999999999, 0, 0, 0, 0, 9281.978999999999),
    @('INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered', 2395.448, 0, 0, 0, 0, 2395.448),
    @('INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible', 0, 0, 0, 0, 0, 0),
    @('JP MORGAN CHASE BANK NA Registered', 124991.729, 0, 0, 0, 0, 124991.729),
    @('JP MORGAN CHASE BANK NA Eligible', 125407.673, 0, 0, 0, 0, 125407.673),
    @('LOOMIS INTERNATIONAL (US) LLC Registered', 68084.33, 0, 0, 0, 0, 68084.33),
    @('LOOMIS INTERNATIONAL (US) LLC Eligible', 106188.481, 0, 0, 0, 0, 106188.481),
    @('MALCA-AMIT USA, LLC Registered', 395.145, 0, 0, 0, 0, 395.145),
    @('MALCA-AMIT USA, LLC Eligible', 0, 0, 0, 0, 0, 0),
    @('MANFRA, TORDELLA & BROOKES, LLC Registered', 54605.27, 0, 0, 0, 0, 54605.27),
    @('MANFRA, TORDELLA & BROOKES, LLC Eligible', 21419.744, 0, 0, 0, 0, 21419.744),
    @('STONEX PRECIOUS METALS LLC Registered', 14122.765, 0, 0, 0, 0, 14122.765),
    @('STONEX PRECIOUS METALS LLC Eligible', 16.075, 0, 0, 0, 0, 16.075)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $entry = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $newDate
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $entry[0]
    $ws.Cells.Item($r, 3).Value = $entry[1]
    $ws.Cells.Item($r, 4).Value = $entry[2]
    $ws.Cells.Item($r, 5).Value = $entry[3]
    $ws.Cells.Item($r, 6).Value = $entry[4]
    $ws.Cells.Item($r, 7).Value = $entry[5]
    $ws.Cells.Item($r, 8).Value = $entry[6]
}

